$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (row 9, algorithm "NB") is dropped entirely in the
# new layout, so remove it first (this also shifts nothing else, since
# it's the very last row).
$ws.Rows.Item(9).Delete()

# --- Header row (B1:H1) ---
# B1:E1 already exist with header text + the bold/centered/bordered
# style; only their text changes. F1:H1 are brand-new header cells that
# need both the new text and that same formatting.
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# Copy the existing header formatting (bold, centered, thin border) onto
# the new header cells without disturbing their (freshly-set) values.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# --- Data rows (A2:H8) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8899968394437421
$ws.Range("D2").Value = 0.02932686691486038
$ws.Range("E2").Value = 0.8872247823860725
$ws.Range("F2").Value = 0.0415375440581336
$ws.Range("G2").Value = 0.872393661384487
$ws.Range("H2").Value = 0.02753491285383908

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8980009481668774
$ws.Range("D3").Value = 0.02737957730481438
$ws.Range("E3").Value = 0.9097542242703532
$ws.Range("F3").Value = 0.04101015707032866
$ws.Range("G3").Value = 0.8732860717264387
$ws.Range("H3").Value = 0.01931091938779279

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8971001896333755
$ws.Range("D4").Value = 0.02404090912375388
$ws.Range("E4").Value = 0.8711469534050179
$ws.Range("F4").Value = 0.0395846454498924
$ws.Range("G4").Value = 0.8778231859883237
$ws.Range("H4").Value = 0.02689385058272304

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7719974715549937
$ws.Range("D5").Value = 0.04768372662983804
$ws.Range("E5").Value = 0.7149257552483359
$ws.Range("F5").Value = 0.04533995245944576
$ws.Range("G5").Value = 0.7411259382819017
$ws.Range("H5").Value = 0.03063760797831392

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8980009481668775
$ws.Range("D6").Value = 0.02465994652500151
$ws.Range("E6").Value = 0.9033282130056325
$ws.Range("F6").Value = 0.04024430996029222
$ws.Range("G6").Value = 0.885137614678899
$ws.Range("H6").Value = 0.02275842567276379

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8554045512010113
$ws.Range("D7").Value = 0.03702121522205542
$ws.Range("E7").Value = 0.8453149001536098
$ws.Range("F7").Value = 0.04973757290997981
$ws.Range("G7").Value = 0.8468557130942452
$ws.Range("H7").Value = 0.02608791639630922

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8873498735777497
$ws.Range("D8").Value = 0.02869691208115329
$ws.Range("E8").Value = 0.8646953405017921
$ws.Range("F8").Value = 0.03974496683500261
$ws.Range("G8").Value = 0.876930775646372
$ws.Range("H8").Value = 0.02204768435220961
